# Update evaluation_metrics.xlsx with the new evaluation results
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.4377224199288256
$wsSummary.Range("C2").Value = 0.08139534883720931
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.1505376344086022
$wsSummary.Range("F2").Value = 0.3070175438596491
$wsSummary.Range("G2").Value = 0.6973180076628352
$wsSummary.Range("H2").Value = 0.7998261102193687
$wsSummary.Range("I2").Value = 28
$wsSummary.Range("J2").Value = 316
$wsSummary.Range("K2").Value = 218
$wsSummary.Range("L2").Value = 0

# --- Sheet 2: Classification Report ---
$wsClassification = $wb.Worksheets.Item("Classification Report")

# Row 2: class "0"
$wsClassification.Range("B2").Value = 1
$wsClassification.Range("C2").Value = 0.4082397003745318
$wsClassification.Range("D2").Value = 0.5797872340425532

# Row 3: class "1"
$wsClassification.Range("B3").Value = 0.08139534883720931
$wsClassification.Range("C3").Value = 1
$wsClassification.Range("D3").Value = 0.1505376344086022

# Row 4: accuracy
$wsClassification.Range("B4").Value = 0.4377224199288256
$wsClassification.Range("C4").Value = 0.4377224199288256
$wsClassification.Range("D4").Value = 0.4377224199288256
$wsClassification.Range("E4").Value = 0.4377224199288256

# Row 5: macro avg
$wsClassification.Range("B5").Value = 0.5406976744186046
$wsClassification.Range("C5").Value = 0.7041198501872659
$wsClassification.Range("D5").Value = 0.3651624342255776

# Row 6: weighted avg
$wsClassification.Range("B6").Value = 0.9542332202267649
$wsClassification.Range("C6").Value = 0.4377224199288256
$wsClassification.Range("D6").Value = 0.5584011329931747

# --- Sheet 3: Confusion Matrix ---
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")

# Row 2: Actual 0
$wsConfusion.Range("B2").Value = 218
$wsConfusion.Range("C2").Value = 316

# Row 3: Actual 1
$wsConfusion.Range("B3").Value = 0
$wsConfusion.Range("C3").Value = 28
